$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Birmingham")

# Trim the stray second address off Emma Melia's e-mail entry so the
# sheet is ready to circulate as a draft.
$ws.Range("E5").Value = "e.melia.1@bham.ac.uk"

# Leave the cursor where the author left it before saving the draft.
$ws.Range("E6").Select()
